# Updates cryptos list prices/volumes (and reorders rows 16-19) per the latest scrape.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cell, $value) {
    # Force the literal text into the cell (avoids Excel auto-coercing
    # numeric-looking strings like "1.003" or "0.000007980" into numbers),
    # then restore the default "Normal" style so no stray number format sticks.
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.Style = "Normal"
}

$ws.Range("D2").Value = "26.653.63"
$ws.Range("E2").Value = "  +3.45%  "
$ws.Range("D3").Value = "1.864.35"
$ws.Range("E3").Value = "  +2.60%  "
Set-TextValue $ws.Range("D4") "1.003"
$ws.Range("E4").Value = "  +0.27%  "
Set-TextValue $ws.Range("D5") "275.30"
$ws.Range("E5").Value = "  -1.46%  "
Set-TextValue $ws.Range("D6") "1.002"
$ws.Range("E6").Value = "  +0.17%  "
Set-TextValue $ws.Range("D7") "0.5276"
$ws.Range("E7").Value = "  +3.45%  "
Set-TextValue $ws.Range("D8") "0.3408"
$ws.Range("E8").Value = "  -3.76%  "
Set-TextValue $ws.Range("D9") "0.06867"
$ws.Range("E9").Value = "  +3.09%  "
Set-TextValue $ws.Range("D10") "19.92"
$ws.Range("E10").Value = "  -0.91%  "
Set-TextValue $ws.Range("D11") "0.7967"
$ws.Range("E11").Value = "  -3.75%  "
Set-TextValue $ws.Range("D12") "0.07731"
$ws.Range("E12").Value = "  -2.19%  "
$ws.Range("D13").Value = "1.892.51"
$ws.Range("E13").Value = "  +4.22%  "
Set-TextValue $ws.Range("D14") "89.69"
$ws.Range("E14").Value = "  +2.17%  "
Set-TextValue $ws.Range("D15") "5.141"
$ws.Range("E15").Value = "  +1.19%  "
$ws.Range("D20").Value = "26.706.00"
$ws.Range("E20").Value = "  +3.49%  "
$ws.Range("D21").Value = "2.111.81"
$ws.Range("E21").Value = "  +4.17%  "
Set-TextValue $ws.Range("D22") "4.733"
$ws.Range("E22").Value = "  -0.43%  "
Set-TextValue $ws.Range("D23") "9.978"
$ws.Range("E23").Value = "  -0.10%  "
Set-TextValue $ws.Range("D24") "6.143"
$ws.Range("E24").Value = "  +0.26%  "
Set-TextValue $ws.Range("D25") "2.350"
$ws.Range("E25").Value = "  +5.48%  "
Set-TextValue $ws.Range("D26") "145.91"
$ws.Range("E26").Value = "  +2.51%  "
Set-TextValue $ws.Range("D27") "17.26"
$ws.Range("E27").Value = "  +0.59%  "
Set-TextValue $ws.Range("D28") "1.650"
$ws.Range("E28").Value = "  -1.16%  "
Set-TextValue $ws.Range("D29") "112.76"
$ws.Range("E29").Value = "  +3.11%  "
Set-TextValue $ws.Range("D30") "4.309"
$ws.Range("E30").Value = "  -0.64%  "
Set-TextValue $ws.Range("D31") "4.314"
$ws.Range("E31").Value = "  +1.80%  "
Set-TextValue $ws.Range("D32") "0.08875"
$ws.Range("E32").Value = "  +1.17%  "
Set-TextValue $ws.Range("D33") "0.04919"
$ws.Range("E33").Value = "  +0.19%  "
Set-TextValue $ws.Range("D34") "1.155"
$ws.Range("E34").Value = "  +1.36%  "
Set-TextValue $ws.Range("D35") "0.7249"
$ws.Range("E35").Value = "  -1.26%  "
Set-TextValue $ws.Range("D36") "2.889"
$ws.Range("E36").Value = "  +0.55%  "
Set-TextValue $ws.Range("D37") "3.243"
$ws.Range("E37").Value = "  +3.21%  "
Set-TextValue $ws.Range("D38") "0.01849"
$ws.Range("E38").Value = "  -0.24%  "
Set-TextValue $ws.Range("D39") "2.312"
$ws.Range("E39").Value = "  -3.18%  "
Set-TextValue $ws.Range("D40") "0.5101"
$ws.Range("E40").Value = "  -1.04%  "
Set-TextValue $ws.Range("D41") "0.9413"
$ws.Range("E41").Value = "  -2.51%  "
Set-TextValue $ws.Range("D42") "115.91"
$ws.Range("E42").Value = "  +4.34%  "
Set-TextValue $ws.Range("D43") "6.126"
$ws.Range("E43").Value = "  -1.83%  "
Set-TextValue $ws.Range("D44") "8.013"
$ws.Range("E44").Value = "  -0.62%  "
$ws.Range("E45").Value = "  +0.14%  "
Set-TextValue $ws.Range("D46") "0.4404"
$ws.Range("E46").Value = "  -3.54%  "
Set-TextValue $ws.Range("D47") "0.1332"
$ws.Range("E47").Value = "  -2.86%  "
Set-TextValue $ws.Range("D48") "9.289"
$ws.Range("E48").Value = "  +0.76%  "
Set-TextValue $ws.Range("D49") "36.18"
$ws.Range("E49").Value = "  -1.29%  "
Set-TextValue $ws.Range("D50") "0.05996"
$ws.Range("E50").Value = "  +3.02%  "
Set-TextValue $ws.Range("D51") "1.477"
$ws.Range("E51").Value = "  -1.79%  "

# Rows 16-19: BinanceUSD/Avalanche and Dai/ShibaInu swapped positions,
# with refreshed prices and volumes.
$ws.Range("B16").Value = "BinanceUSD"
$ws.Range("C16").Value = "https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd"
Set-TextValue $ws.Range("D16") "1.003"
$ws.Range("E16").Value = "  +0.30%  "

$ws.Range("B17").Value = "Avalanche"
$ws.Range("C17").Value = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
Set-TextValue $ws.Range("D17") "14.47"
$ws.Range("E17").Value = "  +2.70%  "

$ws.Range("B18").Value = "Dai"
$ws.Range("C18").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
Set-TextValue $ws.Range("D18") "1.002"
$ws.Range("E18").Value = "  +0.16%  "

$ws.Range("B19").Value = "ShibaInu"
$ws.Range("C19").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
Set-TextValue $ws.Range("D19") "0.000007980"
$ws.Range("E19").Value = "  -0.69%  "

